$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6428.5713
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").Value = $null
$ws.Range("H67").Value = 6428.5713
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").Value = $null
$ws.Range("H132").Value = 2096.55
$ws.Range("I132").Value = 1886.125
$ws.Range("K132").Value = 5658.375
$ws.Range("M132").Value = -3128.375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 995.38464
$ws.Range("I2").Value = 995.38464
$ws.Range("K2").Value = 995.38464
$ws.Range("M2").Value = -882.38464
$ws.Range("H32").Value = 2340.7476
$ws.Range("I32").Value = 2340.7476
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2340.7476
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = $null
$ws.Range("N32").Value = -2053.7476
$ws.Range("H45").Value = 2249
$ws.Range("I45").Value = 2748
$ws.Range("K45").Value = 2748
$ws.Range("M45").Value = -2371
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = $null
$ws.Range("N103").Value = 0
$ws.Range("H116").Value = 995.38464
$ws.Range("I116").Value = 995.38464
$ws.Range("K116").Value = 995.38464
$ws.Range("M116").Value = 1298.61536
$ws.Range("H122").Value = 4711.2915
$ws.Range("I122").Value = 3283.7273
$ws.Range("J122").Value = 5919.231
$ws.Range("K122").Value = 9851.1819
$ws.Range("L122").Value = 17757.693
$ws.Range("M122").Value = -7401.1819
$ws.Range("N122").Value = -22657.693
$ws.Range("H133").Value = 54916.5
$ws.Range("J133").Value = 54916.5
$ws.Range("L133").Value = 54916.5
$ws.Range("N133").Value = -59976.5
$ws.Range("H137").Value = 88750
$ws.Range("J137").Value = 88750
$ws.Range("L137").Value = 88750
$ws.Range("N137").Value = -98950

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 995.38464
$ws.Range("I3").Value = 995.38464
$ws.Range("K3").Value = 995.38464
$ws.Range("M3").Value = -881.38464
$ws.Range("H107").Value = 772136.3
$ws.Range("I107").Value = 2523.5454
$ws.Range("K107").Value = 2523.5454
$ws.Range("M107").Value = -603.5454
$ws.Range("H134").Value = 92201.25
$ws.Range("I134").Value = 9641.5
$ws.Range("J134").Value = 505000
$ws.Range("K134").Value = 28924.5
$ws.Range("L134").Value = 1515000
$ws.Range("M134").Value = -26389.5
$ws.Range("N134").Value = -1520070

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4035.6
$ws.Range("I99").Value = 3459.6667
$ws.Range("J99").Value = 4899.5
$ws.Range("K99").Value = 3459.6667
$ws.Range("L99").Value = 4899.5
$ws.Range("M99").Value = -1961.6667
$ws.Range("N99").Value = -7895.5
$ws.Range("H126").Value = 4035.6
$ws.Range("I126").Value = 3459.6667
$ws.Range("J126").Value = 4899.5
$ws.Range("K126").Value = 10379.0001
$ws.Range("L126").Value = 14698.5
$ws.Range("M126").Value = -7909.000100000001
$ws.Range("N126").Value = -19638.5
$ws.Range("H133").Value = 47538.23
$ws.Range("I133").Value = 39999.5
$ws.Range("J133").Value = 48908.91
$ws.Range("K133").Value = 39999.5
$ws.Range("L133").Value = 48908.91
$ws.Range("M133").Value = -37469.5
$ws.Range("N133").Value = -53968.91
$ws.Range("H134").Value = 419964.03
$ws.Range("I134").Value = 3455.6191
$ws.Range("K134").Value = 10366.8573
$ws.Range("M134").Value = -7831.8573

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2399.8
$ws.Range("I80").Value = 1501
$ws.Range("K80").Value = 4503
$ws.Range("M80").Value = -3567
$ws.Range("H83").Value = 2399.8
$ws.Range("I83").Value = 1501
$ws.Range("K83").Value = 13509
$ws.Range("M83").Value = -8829
$ws.Range("H87").Value = 18999.666
$ws.Range("I87").Value = 18999.666
$ws.Range("K87").Value = 56998.99800000001
$ws.Range("M87").Value = -55750.99800000001
$ws.Range("H90").Value = 18999.666
$ws.Range("I90").Value = 18999.666
$ws.Range("K90").Value = 170996.994
$ws.Range("M90").Value = -164756.994
$ws.Range("H93").Value = 340
$ws.Range("I93").Value = 340
$ws.Range("K93").Value = 1020
$ws.Range("M93").Value = 852
$ws.Range("H114").Value = 738.5
$ws.Range("I114").Value = 332
$ws.Range("J114").Value = 912.7143
$ws.Range("K114").Value = 996
$ws.Range("L114").Value = 2738.1429
$ws.Range("M114").Value = 2258
$ws.Range("N114").Value = -9246.142899999999
$ws.Range("H131").Value = 3121.6304
$ws.Range("I131").Value = 2069.875
$ws.Range("J131").Value = 3343.0527
$ws.Range("K131").Value = 6209.625
$ws.Range("L131").Value = 10029.1581
$ws.Range("M131").Value = -1169.625
$ws.Range("N131").Value = -20109.1581

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8383.950000000001
$ws.Range("I70").Value = 6690.923
$ws.Range("K70").Value = 6690.923
$ws.Range("M70").Value = -6420.923
$ws.Range("H73").Value = 8383.950000000001
$ws.Range("I73").Value = 6690.923
$ws.Range("K73").Value = 6690.923
$ws.Range("M73").Value = -5754.923
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = $null
$ws.Range("N95").Value = 0
$ws.Range("H102").Value = 4657.6
$ws.Range("I102").Value = 2462.7144
$ws.Range("J102").Value = 6578.125
$ws.Range("K102").Value = 2462.7144
$ws.Range("L102").Value = 6578.125
$ws.Range("M102").Value = -840.7143999999998
$ws.Range("N102").Value = -9822.125
$ws.Range("H132").Value = 52035.57
$ws.Range("I132").Value = 4707.5557
$ws.Range("J132").Value = 336003.66
$ws.Range("K132").Value = 14122.6671
$ws.Range("L132").Value = 1008010.98
$ws.Range("M132").Value = -11592.6671
$ws.Range("N132").Value = -1013070.98

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6926.2144
$ws.Range("I7").Value = 7645.5
$ws.Range("J7").Value = 6386.75
$ws.Range("K7").Value = 7645.5
$ws.Range("L7").Value = 6386.75
$ws.Range("M7").Value = -7533.5
$ws.Range("N7").Value = -6610.75
$ws.Range("H40").Value = 4428.355
$ws.Range("I40").Value = 3713.9583
$ws.Range("J40").Value = 6877.7144
$ws.Range("K40").Value = 3713.9583
$ws.Range("L40").Value = 6877.7144
$ws.Range("M40").Value = -3577.9583
$ws.Range("N40").Value = -7149.7144
$ws.Range("H55").Value = 2140.875
$ws.Range("I55").Value = 69.5
$ws.Range("K55").Value = 69.5
$ws.Range("M55").Value = 103.5
$ws.Range("H61").Value = 3355.3333
$ws.Range("I61").Value = 3355.3333
$ws.Range("K61").Value = 3355.3333
$ws.Range("M61").Value = -3153.3333
$ws.Range("H113").Value = 3355.3333
$ws.Range("I113").Value = 3355.3333
$ws.Range("K113").Value = 3355.3333
$ws.Range("M113").Value = -1185.3333
$ws.Range("H126").Value = 6926.2144
$ws.Range("I126").Value = 7645.5
$ws.Range("J126").Value = 6386.75
$ws.Range("K126").Value = 22936.5
$ws.Range("L126").Value = 19160.25
$ws.Range("M126").Value = -20466.5
$ws.Range("N126").Value = -24100.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6141.643
$ws.Range("I62").Value = 5097.2
$ws.Range("J62").Value = 6721.8887
$ws.Range("K62").Value = 5097.2
$ws.Range("L62").Value = 6721.8887
$ws.Range("M62").Value = -4473.2
$ws.Range("N62").Value = -7969.8887
$ws.Range("H65").Value = 6141.643
$ws.Range("I65").Value = 5097.2
$ws.Range("J65").Value = 6721.8887
$ws.Range("K65").Value = 25486
$ws.Range("L65").Value = 33609.4435
$ws.Range("M65").Value = -22366
$ws.Range("N65").Value = -39849.4435
$ws.Range("H122").Value = 30306838
$ws.Range("I122").Value = 38465204
$ws.Range("K122").Value = 115395612
$ws.Range("M122").Value = -115393162
$ws.Range("H126").Value = 999.5
$ws.Range("I126").Value = 999.5
$ws.Range("K126").Value = 2998.5
$ws.Range("M126").Value = -528.5
$ws.Range("H136").Value = 114132.164
$ws.Range("I136").Value = 3198.7646
$ws.Range("K136").Value = 9596.293799999999
$ws.Range("M136").Value = -7046.293799999999
